# This workbook's data rows (2-20) got reshuffled between source rows and
# target rows: each target row's columns D..Q (Fecha, Codreg, Categoria ID,
# Categoria, Variedad, Calidad, Volumen, Precio minimo, Precio maximo,
# Precio promedio ponderado, Unidad de comercializacion, Origen, Precio $/Kg,
# Kg o Unidades) end up holding the values that used to live in a different
# row. Columns A, B, C and R (Mercado ID, Mercado, Region, Clasificacion) do
# not change. Row 16 is untouched.
#
# Mapping: target row -> source row (the row whose original D..Q values now
# appear in the target row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$mapping = @{
    2  = 9
    3  = 12
    4  = 8
    5  = 10
    6  = 11
    7  = 6
    8  = 19
    9  = 17
    10 = 18
    11 = 20
    12 = 5
    13 = 7
    14 = 13
    15 = 3
    16 = 16
    17 = 2
    18 = 14
    19 = 4
    20 = 15
}

# Columns D (4) through Q (17) hold the data that gets shuffled around.
$firstCol = 4
$lastCol = 17

# Snapshot the original values for every row before writing anything, since
# several rows read from each other (a straightforward row-by-row rewrite
# would clobber data that a later row still needs to read).
$snapshot = @{}
for ($row = 2; $row -le 20; $row++) {
    $rowValues = @{}
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $rowValues[$col] = $ws.Cells.Item($row, $col).Value2
    }
    $snapshot[$row] = $rowValues
}

foreach ($targetRow in $mapping.Keys) {
    $sourceRow = $mapping[$targetRow]
    if ($sourceRow -eq $targetRow) {
        continue
    }
    $sourceValues = $snapshot[$sourceRow]
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $ws.Cells.Item($targetRow, $col).Value = $sourceValues[$col]
    }
}
